# MRD-1840 Add emails for revocation order to Part A
#
# The paragraph that used to hold a FORMTEXT form-field for the
# "Email address for receipt of a copy of the revocation order to be
# sent to:" question is converted to a merge-field placeholder
# ({{revocation_order_recipients}}) on a new line, and the paragraph
# mark becomes bold.

$d = $word.ActiveDocument

# Locate the target paragraph robustly by its (unique) label text
# rather than a hard-coded paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Email address for receipt of a copy of the revocation order to be sent to:*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'Email address for receipt of a copy of the revocation order' paragraph"
}

# Replace the paragraph's contents (label run kept, everything after it
# replaced) with: the bold paragraph mark, a line break, and the
# {{revocation_order_recipients}} merge field (the placeholder word is
# wrapped in proofErr spell-check markers, matching the convention used
# elsewhere in this document).
$newParaXml = '<w:p w14:paraId="765FE164" w14:textId="77777777" w:rsidR="00AC158B" w:rsidRDefault="00AC158B" w:rsidP="00F23ECD">' + `
  '<w:pPr>' + `
    '<w:pBdr>' + `
      '<w:top w:val="single" w:sz="4" w:space="1" w:color="auto"/>' + `
      '<w:left w:val="single" w:sz="4" w:space="4" w:color="auto"/>' + `
      '<w:bottom w:val="single" w:sz="4" w:space="1" w:color="auto"/>' + `
      '<w:right w:val="single" w:sz="4" w:space="4" w:color="auto"/>' + `
    '</w:pBdr>' + `
    '<w:tabs><w:tab w:val="left" w:pos="1496"/></w:tabs>' + `
    '<w:spacing w:line="360" w:lineRule="auto"/>' + `
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r w:rsidRPr="00823C42">' + `
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:color w:val="800080"/></w:rPr>' + `
    '<w:t>Email address for receipt of a copy of the revocation order to be sent to:</w:t>' + `
  '</w:r>' + `
  '<w:r w:rsidRPr="005D7929">' + `
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/></w:rPr>' + `
    '<w:br/>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' + `
    '<w:t>{{</w:t>' + `
  '</w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r>' + `
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' + `
    '<w:t>revocation_order_recipients</w:t>' + `
  '</w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + `
    '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr>' + `
    '<w:t>}}</w:t>' + `
  '</w:r>' + `
'</w:p>'

$target.Range.InsertXML($newParaXml)
